$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.123.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4666"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2815"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06418"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.25"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "96.65"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +13.69%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07550"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.841.81"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6362"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "293.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +20.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.124.08"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007363"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.098.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.053"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.82"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.116"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.921"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1086"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.330"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.006"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.801"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04918"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7243"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.110"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.738"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01919"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.663"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8614"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.961"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.09"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.611"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4045"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.049"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.941"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1186"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "835.46"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +17.08%  "
